$d = $word.ActiveDocument

# Locate the paragraph that holds the "Ver no Jupiter ..." text and the one
# right after it holding the "(C) 2020 ..." footer text. Together with the
# blank paragraph that precedes the first of them, these three paragraphs
# are being dropped from the page (the site's generated footer block).
$count = $d.Paragraphs.Count
$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Ver no Jupiter*") {
        $startIndex = $i - 1
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $endIndex = $i
    }
}

if ($startIndex -ge 1 -and $endIndex -ge $startIndex) {
    $rangeStart = $d.Paragraphs.Item($startIndex).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIndex).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
